$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-add the block that was previously removed: a duplicate of rows 1-13
# (with the same blank-row gaps) starting at row 14, leaving row 27 blank
# (mirroring the old trailing blank row 14).
$pairs = @(
    @("A1", "A14"),
    @("A2", "A15"),
    @("A3", "A16"),
    @("A5", "A18"),
    @("A6", "A19"),
    @("A7", "A20"),
    @("A8", "A21"),
    @("B10", "B23"),
    @("D10", "D23"),
    @("E10", "E23"),
    @("F10", "F23"),
    @("B11", "B24"),
    @("D11", "D24"),
    @("E11", "E24"),
    @("F11", "F24"),
    @("A13", "A26")
)

foreach ($pair in $pairs) {
    $srcAddr = $pair[0]
    $dstAddr = $pair[1]

    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($dstAddr)

    # Force text interpretation first so values that look numeric/date-like
    # (e.g. "1/13/2020") land as literal text, same as the source cell,
    # instead of being auto-coerced into a date serial number.
    $dst.NumberFormat = "@"
    $dst.Value = $src.Text

    # Now copy the real formatting (style) from the source cell over the
    # destination without disturbing any other cells/rows. Because the
    # destination is already stored as text, reapplying the source's
    # number format here does not convert it back into a number.
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
